$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Val)
    $origStyle = $Cell.Style
    $Cell.NumberFormat = "@"
    $Cell.Value = $Val
    $Cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "25.900.77"
Set-TextValue $ws.Range("E2") "  -0.68%  "
Set-TextValue $ws.Range("D3") "1.740.61"
Set-TextValue $ws.Range("E3") "  -0.47%  "
Set-TextValue $ws.Range("D4") "1.002"
Set-TextValue $ws.Range("E4") "  +0.32%  "
Set-TextValue $ws.Range("D5") "249.39"
Set-TextValue $ws.Range("E5") "  +6.32%  "
Set-TextValue $ws.Range("D6") "1.001"
Set-TextValue $ws.Range("D7") "0.5101"
Set-TextValue $ws.Range("E7") "  -3.87%  "
Set-TextValue $ws.Range("D8") "0.2742"
Set-TextValue $ws.Range("E8") "  -1.72%  "
Set-TextValue $ws.Range("D9") "0.06180"
Set-TextValue $ws.Range("E9") "  -0.18%  "
Set-TextValue $ws.Range("D10") "1.742.79"
Set-TextValue $ws.Range("E10") "  -0.29%  "
Set-TextValue $ws.Range("D11") "0.07225"
Set-TextValue $ws.Range("E11") "  -0.05%  "
Set-TextValue $ws.Range("D12") "15.08"
Set-TextValue $ws.Range("E12") "  -2.07%  "
Set-TextValue $ws.Range("D13") "0.6478"
Set-TextValue $ws.Range("E13") "  +0.36%  "
Set-TextValue $ws.Range("D14") "4.623"
Set-TextValue $ws.Range("E14") "  +0.15%  "
Set-TextValue $ws.Range("D15") "77.47"
Set-TextValue $ws.Range("E15") "  -1.23%  "
Set-TextValue $ws.Range("D16") "1.002"
Set-TextValue $ws.Range("E16") "  +0.19%  "
Set-TextValue $ws.Range("D17") "1.002"
Set-TextValue $ws.Range("E17") "  +0.31%  "
Set-TextValue $ws.Range("D18") "25.935.91"
Set-TextValue $ws.Range("E18") "  -0.17%  "
Set-TextValue $ws.Range("D19") "11.81"
Set-TextValue $ws.Range("E19") "  +1.54%  "
Set-TextValue $ws.Range("D20") "0.000006797"
Set-TextValue $ws.Range("E20") "  +0.81%  "
Set-TextValue $ws.Range("D21") "1.968.32"
Set-TextValue $ws.Range("E21") "  -0.01%  "
Set-TextValue $ws.Range("D22") "4.263"
Set-TextValue $ws.Range("E22") "  -1.27%  "
Set-TextValue $ws.Range("D23") "8.662"
Set-TextValue $ws.Range("E23") "  -1.46%  "
Set-TextValue $ws.Range("D24") "5.382"
Set-TextValue $ws.Range("E24") "  +3.20%  "
Set-TextValue $ws.Range("D25") "135.58"
Set-TextValue $ws.Range("E25") "  -2.34%  "
Set-TextValue $ws.Range("D26") "1.502"
Set-TextValue $ws.Range("E26") "  -0.55%  "
Set-TextValue $ws.Range("E27") "  -0.41%  "
Set-TextValue $ws.Range("D28") "1.774"
Set-TextValue $ws.Range("E28") "  -1.90%  "
Set-TextValue $ws.Range("D29") "105.92"
Set-TextValue $ws.Range("E29") "  +1.07%  "
Set-TextValue $ws.Range("D30") "3.927"
Set-TextValue $ws.Range("E30") "  +3.37%  "
Set-TextValue $ws.Range("D31") "0.08214"
Set-TextValue $ws.Range("E31") "  -1.16%  "
Set-TextValue $ws.Range("E32") "  -1.28%  "
Set-TextValue $ws.Range("D33") "0.04699"
Set-TextValue $ws.Range("D34") "2.656"
Set-TextValue $ws.Range("E34") "  +0.35%  "
Set-TextValue $ws.Range("D35") "0.9963"
Set-TextValue $ws.Range("E35") "  -0.55%  "
Set-TextValue $ws.Range("D36") "0.6248"
Set-TextValue $ws.Range("E36") "  -1.44%  "
Set-TextValue $ws.Range("D37") "2.739"
Set-TextValue $ws.Range("E37") "  +1.06%  "
Set-TextValue $ws.Range("D38") "0.01608"
Set-TextValue $ws.Range("E38") "  +0.84%  "
Set-TextValue $ws.Range("D39") "1.913"
Set-TextValue $ws.Range("E39") "  -1.58%  "
Set-TextValue $ws.Range("D40") "1.001"
Set-TextValue $ws.Range("E40") "  +0.16%  "
Set-TextValue $ws.Range("D41") "99.82"
Set-TextValue $ws.Range("E41") "  +1.06%  "
Set-TextValue $ws.Range("D42") "0.7551"
Set-TextValue $ws.Range("E42") "  +1.97%  "
Set-TextValue $ws.Range("D43") "0.3837"
Set-TextValue $ws.Range("E43") "  -2.07%  "
Set-TextValue $ws.Range("D44") "4.999"
Set-TextValue $ws.Range("E44") "  -0.54%  "
Set-TextValue $ws.Range("B45") "Algorand"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D45") "0.1128"
Set-TextValue $ws.Range("E45") "  -1.71%  "
Set-TextValue $ws.Range("B46") "Aptos"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D46") "6.290"
Set-TextValue $ws.Range("E46") "  -0.57%  "
Set-TextValue $ws.Range("D47") "55.37"
Set-TextValue $ws.Range("E47") "  +2.30%  "
Set-TextValue $ws.Range("D48") "0.05228"
Set-TextValue $ws.Range("E48") "  -2.18%  "
Set-TextValue $ws.Range("D49") "30.78"
Set-TextValue $ws.Range("E49") "  +0.12%  "
Set-TextValue $ws.Range("D50") "7.516"
Set-TextValue $ws.Range("E50") "  -1.81%  "
Set-TextValue $ws.Range("D51") "0.3412"
Set-TextValue $ws.Range("E51") "  -1.45%  "
